$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Notes cell for the MCU row to mention the Segger J-Link EDU programmer
$ws.Range("C4").Value = "Includes FRDM-K22 dev kit and Segger J-Link EDU"

# Add the Programmer cost (70) to the MCU subsystem cost formula
$ws.Range("B4").Formula = "=10+5+30+70"

# Update the active selection to C5
$ws.Range("C5").Select()
